$d = $word.ActiveDocument

# Insert the new "Novembre" section (title, sub-title and three new bullet
# points) right at the very end of the document body, after the existing
# last paragraph ("Suite de la création des fonctions ...").
#
# Using Range.InsertXML on a collapsed range positioned at the exact end of
# the document content lets us add the five new paragraphs with their
# final pPr/style/numbering already correct in one shot, without touching
# (or leaving stray formatting marks on) the paragraph that precedes them.

$end = $d.Content.End
$insertionPoint = $d.Range($end, $end)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newParagraphsXml =
  "<w:p $ns><w:pPr><w:pStyle w:val=`"Titre1`"/></w:pPr><w:r><w:t>Novembre</w:t></w:r></w:p>" +
  "<w:p $ns><w:pPr><w:pStyle w:val=`"Titre2`"/></w:pPr><w:r><w:t>9 Novembre 17</w:t></w:r></w:p>" +
  "<w:p $ns><w:pPr><w:pStyle w:val=`"Paragraphedeliste`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"2`"/></w:numPr></w:pPr><w:r><w:t>RAZ du projet</w:t></w:r></w:p>" +
  "<w:p $ns><w:pPr><w:pStyle w:val=`"Paragraphedeliste`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"2`"/></w:numPr></w:pPr><w:r><w:t>Création de la page de résultat</w:t></w:r></w:p>" +
  "<w:p $ns><w:pPr><w:pStyle w:val=`"Paragraphedeliste`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"2`"/></w:numPr></w:pPr><w:r><w:t>Adaptation de la page de validation des frais</w:t></w:r></w:p>"

$insertionPoint.InsertXML($newParagraphsXml)

$d.Save()
